$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.121.33'
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").Value = '2.588.02'
$ws.Range("E3").Value = '  +2.92%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.30%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.23%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("D13").Value = '2.987.48'
$ws.Range("E13").Value = '  +2.99%  '

$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("D15").Value = '2.520.36'
$ws.Range("E15").Value = '  +1.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '43.210.27'
$ws.Range("E18").Value = '  +1.31%  '

$ws.Range("E19").Value = '  +2.28%  '

$ws.Range("E20").Value = '  -1.74%  '

$ws.Range("D21").Value = '0.0₃0966'
$ws.Range("E21").Value = '  +1.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  +3.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.07%  '

$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("E28").Value = '  -0.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.79'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  +0.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("E35").Value = '  +3.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.22%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.113'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.15%  '

$ws.Range("B39").Value = 'ApeXProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.60%  '

$ws.Range("E40").Value = '  +0.39%  '

$ws.Range("E41").Value = '  -1.71%  '

$ws.Range("E42").Value = '  +5.33%  '

$ws.Range("E43").Value = '  +0.91%  '

$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("D45").Value = '2.012.55'
$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.26%  '

$ws.Range("E47").Value = '  +1.40%  '

$ws.Range("D48").Value = '2.838.36'
$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.36'
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  +4.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.11%  '
